# Append the 11/21/2025 Kaspa buy row (row 15) to the log sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds the date as plain text (matching the other rows in this
# sheet), not an Excel date serial. Temporarily force a text number format
# so Excel doesn't auto-parse the "mm/dd/yyyy" string into a date, then
# clear the formatting again so the new cell stays styleless like its
# siblings (A2:A5, A7:A14).
$ws.Cells.Item(15, 1).NumberFormat = "@"
$ws.Cells.Item(15, 1).Value = "11/21/2025"
$ws.Cells.Item(15, 1).ClearFormats()

$ws.Cells.Item(15, 2).Value = 594.4650000000001
$ws.Cells.Item(15, 3).Value = 0.04163407433574726
$ws.Cells.Item(15, 4).Value = 25
